# Quarterly database update: drop the oldest quarter ("فصل سوم منتهی به
# 1399/06"), shift every remaining quarter's data one column to the left,
# and append the newest quarter ("فصل اول منتهی به 1401/12") with its data
# at the end (column N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quarter header labels (rows 8 and 24, columns E..N) -------------------
$quarters = @(
    "فصل چهارم منتهی به 1399/09",
    "فصل اول منتهی به 1399/12",
    "فصل دوم منتهی به 1400/03",
    "فصل سوم منتهی به 1400/06",
    "فصل چهارم منتهی به 1400/09",
    "فصل اول منتهی به 1400/12",
    "فصل دوم منتهی به 1401/03",
    "فصل سوم منتهی به 1401/06",
    "فصل چهارم منتهی به 1401/09",
    "فصل اول منتهی به 1401/12"
)

for ($i = 0; $i -lt $quarters.Length; $i++) {
    $col = 5 + $i  # E=5 .. N=14
    $ws.Cells.Item(8, $col).Value = $quarters[$i]
    $ws.Cells.Item(24, $col).Value = $quarters[$i]
}

# --- Data rows, shifted left by one quarter with the new quarter appended --
$dataRows = @{
    10 = @(4656, 0, 118255, 20121, 40559, 10444, 75739, 0, 187686, 22254)
    11 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    12 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    13 = @(10797, 0, 3863, 6185, 14722, 7389, 8515, 10199, 12283, 8500)
    14 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    15 = @(2391, 0, 1297, 0, 2475, 214, 1665, 16, 2930, 1222)
    16 = @(621, 539, 705, 404, 1173, 407, 1298, 59, 1504, 941)
    17 = @(24962, 22625, 27079, 25789, 29322, 28041, 40673, 31856, 59204, 46733)
    18 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    19 = @(115537, 42272, 17841, 28002, 57582, 54011, 49548, 109058, 46135, 55609)
    20 = @(158964, 65436, 169040, 80501, 145833, 100506, 177438, 151188, 309742, 135259)
    26 = @(599, 603, 616, 618, 614, 584, 834, 834, 568, 569)
    27 = @(99, 96, 96, 96, 96, 101, 99, 99, 100, 94)
}

foreach ($r in $dataRows.Keys) {
    $values = $dataRows[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i  # E=5 .. N=14
        $ws.Cells.Item([int]$r, $col).Value = $values[$i]
    }
}
